# Add files via upload — add Ryan Huang's two papers to the "Main" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# ------------------------------------------------------------------
# Row 19 — Coenzyme Q10 / Leber's hereditary optic neuropathy paper
# ------------------------------------------------------------------
$ws.Range("B19").Value = "Coenzyme Q10 trapping in mitochondrial complex I underlies Leber's hereditary optic neuropathy"
$ws.Range("B19").Borders.LineStyle = -4142   # xlLineStyleNone -> drop the grid border on B19

$ws.Range("C19").Value = 2023

$ws.Range("D19").Value = "https://doi.org/10.1073/pnas.2304884120"
$ws.Range("D19").WrapText = $true
$ws.Range("D19").Borders.LineStyle = -4142
$ws.Hyperlinks.Add($ws.Range("D19"), "https://doi.org/10.1073/pnas.2304884120") | Out-Null

$ws.Range("E19").Value = "Ryan Huang"

$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = "N/A"
$ws.Range("J19").Value = 5
$ws.Range("K19").Value = "N/A (?)"
$ws.Range("L19").Value = 5
$ws.Range("M19").Value = "This is a molecular dynamics simulation study that publishes all its data and codes in an online repository. However, the actual steps performed in the simulation were not described in detail in the paper."

$ws.Rows.Item(19).RowHeight = 33

# ------------------------------------------------------------------
# Row 20 — morphology of an infarct in NAION paper
# ------------------------------------------------------------------
$ws.Range("B20").Value = "The morphology of an infarct in nonarteritic anterior ischemic optic neuropathy"
$ws.Range("C20").Value = 2003

$ws.Range("D20").Value = "https://doi.org/10.1016/S0161-6420(03)00804-2"
$ws.Range("D20").WrapText = $true
$ws.Range("D20").Borders.LineStyle = -4142
$ws.Hyperlinks.Add($ws.Range("D20"), "https://doi.org/10.1016/S0161-6420(03)00804-2") | Out-Null

$ws.Range("E20").Value = "Ryan Huang"

$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 5
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 3
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 3
$ws.Range("M20").Value = "This is a pathology examination study which is based on only one patient. It is difficult to replicate since it is very hard to obtain  such patient samples. However, the pathological examination protocol was well-described in the paper."

$ws.Rows.Item(20).RowHeight = 36

# ------------------------------------------------------------------
# Rows 21-23 — fill in student name only (rest left blank)
# ------------------------------------------------------------------
$ws.Range("E21").Value = "Ryan Huang"
$ws.Range("E22").Value = "Ryan Huang"
$ws.Range("E23").Value = "Ryan Huang"

# ------------------------------------------------------------------
# Misc formatting touch-ups
# ------------------------------------------------------------------
$ws.Rows.Item(18).RowHeight = 16
$ws.Columns.Item(4).ColumnWidth = 11.67

$ws.Range("M21").Select()
